# Workbook/worksheet handles (workbook is already open per harness contract)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update threshold values ---
$ws.Range("B2").Value = 5
$ws.Range("B3").Value = 4.5
$ws.Range("C4").Value = 1.5

# Remove the "theta_threshold_range" row (row 5) entirely -- this both
# drops that now-unused shared string from the workbook and shifts the
# former row 6 ("pie_threshold_range") up into row 5.
$ws.Range("A5").EntireRow.Delete()

# The shifted-up row (now row 5) gets its Min/Max updated.
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 20

# --- Re-fit the remaining columns to their (now shorter) content ---
$ws.Columns("A:C").AutoFit()

# Leave the selection on the last edited cell, matching where the author
# finished working.
$ws.Range("C5").Select() | Out-Null
